$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh country statistics (new case counts as of 12 Oct 2020, 13:00)
# Estados Unidos
$ws.Range("B4").Value = 7992932
$ws.Range("C4").Value = 934
$ws.Range("D4").Value = 5128492
$ws.Range("E4").Value = 2644734
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 219706

# India
$ws.Range("B5").Value = 7122862
$ws.Range("C5").Value = 3562
$ws.Range("E5").Value = 864118
$ws.Range("G5").Value = 25
$ws.Range("H5").Value = 109209

# Iran
$ws.Range("B16").Value = 504281
$ws.Range("C16").Value = 4206
$ws.Range("D16").Value = 409121
$ws.Range("E16").Value = 66344
$ws.Range("G16").Value = 272
$ws.Range("H16").Value = 28816

# Banglades
$ws.Range("B19").Value = 379738
$ws.Range("C19").Value = 1472
$ws.Range("D19").Value = 294391
$ws.Range("E19").Value = 79792
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 5555

# Alemania
$ws.Range("D25").Value = 276900
$ws.Range("E25").Value = 39880

# Rumania
$ws.Range("B32").Value = 157352
$ws.Range("C32").Value = 2069
$ws.Range("D32").Value = 118912
$ws.Range("E32").Value = 32973
$ws.Range("G32").Value = 56
$ws.Range("H32").Value = 5467

# Catar
$ws.Range("B37").Value = 128191
$ws.Range("C37").Value = 206
$ws.Range("D37").Value = 125176
$ws.Range("E37").Value = 2795

# Emiratos Arabes Unidos
$ws.Range("B45").Value = 107293
$ws.Range("C45").Value = 1064
$ws.Range("D45").Value = 98555
$ws.Range("E45").Value = 8292
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 446

# Suiza
$ws.Range("B60").Value = 64436
$ws.Range("C60").Value = 4068
$ws.Range("E60").Value = 13945
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 2091

# Libia
$ws.Range("B71").Value = 43821
$ws.Range("C71").Value = 1109
$ws.Range("D71").Value = 24466
$ws.Range("E71").Value = 18711
$ws.Range("G71").Value = 13
$ws.Range("H71").Value = 644

# Malasia
$ws.Range("B94").Value = 16220
$ws.Range("C94").Value = 563
$ws.Range("D94").Value = 11022
$ws.Range("E94").Value = 5039
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 159

# Senegal
$ws.Range("B98").Value = 15292
$ws.Range("C98").Value = 24
$ws.Range("D98").Value = 13390
$ws.Range("E98").Value = 1587
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 315

# Georgia
$ws.Range("E101").Value = 5641
$ws.Range("G101").Value = 8
$ws.Range("H101").Value = 93

# Finlandia
$ws.Range("B102").Value = 12212
$ws.Range("C102").Value = 214
$ws.Range("E102").Value = 3366

# Eslovenia
$ws.Range("B114").Value = 8832
$ws.Range("C114").Value = 169
$ws.Range("D114").Value = 5219
$ws.Range("E114").Value = 3444
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 169

# Malta
$ws.Range("B142").Value = 3844
$ws.Range("C142").Value = 68
$ws.Range("D142").Value = 2981
$ws.Range("E142").Value = 820
$ws.Range("G142").Value = 2
$ws.Range("H142").Value = 43

# Gibraltar
$ws.Range("B180").Value = 486
$ws.Range("C180").Value = 1
$ws.Range("D180").Value = 419
$ws.Range("E180").Value = 67

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 13:00"

# Re-sort the country table by Casos totales (column B) descending, as in the live feed
$sortRange = $ws.Range("A4:H220")
$keyRange = $ws.Range("B4:B220")
$sortRange.Sort($keyRange, 2)
